$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells in column D whose new value would otherwise be auto-detected
# by Excel as a number (single decimal point) need to be forced to stay text,
# matching the source data which stores prices as plain text strings.
$textDCells = @("D5","D8","D14","D16","D18","D19","D22","D24","D26","D28","D38","D39","D43","D44","D47","D50")
foreach ($addr in $textDCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.039.12"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.562.87"
$ws.Range("E3").Value = "  +0.46%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.36%  "

# Row 5 - BNB
$ws.Range("D5").Value = "208.49"
$ws.Range("E5").Value = "  +0.66%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.55%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.42%  "

# Row 8 - Solana
$ws.Range("D8").Value = "22.08"
$ws.Range("E8").Value = "  -0.21%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.91%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.82%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.19%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.555.90"
$ws.Range("E12").Value = "  +0.02%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  -0.25%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "0.521"
$ws.Range("E14").Value = "  +0.20%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "27.040.02"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "61.93"

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +1.62%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "215.75"
$ws.Range("E18").Value = "  -1.08%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "7.38"
$ws.Range("E19").Value = "  +1.02%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.37%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +1.89%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "9.21"
$ws.Range("E22").Value = "  -0.17%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  +0.04%  "

# Row 24 - Monero
$ws.Range("D24").Value = "153.24"

# Row 25 - Cosmos
$ws.Range("E25").Value = "  -0.58%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "15.05"
$ws.Range("E26").Value = "  +0.73%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +1.39%  "

# Row 28 - BinanceUSD
$ws.Range("D28").Value = "1.01"
$ws.Range("E28").Value = "  +0.35%  "

# Row 29 - Hedera
$ws.Range("E29").Value = "  +1.11%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +3.09%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  +0.00%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +3.33%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.434.11"
$ws.Range("E33").Value = "  +0.45%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +1.33%  "

# Row 35 - TrustWalletToken
$ws.Range("E35").Value = "  +8.49%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +2.54%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +0.93%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "0.534"
$ws.Range("E38").Value = "  +2.15%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "5.91"
$ws.Range("E39").Value = "  +1.94%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -0.17%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.46%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  +0.88%  "

# Row 43 - MXToken
$ws.Range("D43").Value = "2.32"
$ws.Range("E43").Value = "  -0.10%  "

# Row 44 - Aave
$ws.Range("D44").Value = "64.75"
$ws.Range("E44").Value = "  +0.55%  "

# Row 45 - RenderToken
$ws.Range("E45").Value = "  -0.66%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "1.698.57"
$ws.Range("E46").Value = "  +0.42%  "

# Row 47 - Quant
$ws.Range("D47").Value = "86.96"
$ws.Range("E47").Value = "  -1.01%  "

# Row 48 - BabyDogeCoin
$ws.Range("E48").Value = "  +3.78%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  -0.72%  "

# Row 50 - Algorand
$ws.Range("D50").Value = "0.0960"
$ws.Range("E50").Value = "  +0.38%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  +0.37%  "
